# Update "output/StructureDefinition-social-context.xlsx" for the 2025-08 refresh:
#   - new canonical base URL (2rdoc.pt IG instead of the old GitHub shorthand URL)
#   - refreshed publication Date
#   - new ValueSet canonical URL
#   - column widths on the "Elements" sheet were re-measured by the IG
#     publisher (bestFit) after the content changed, so they shrink slightly
#
# Sheet "Metadata" (sheet1) holds the Property/Value table (URL in B2, Date in B8).
# Sheet "Elements" (sheet2) holds the per-element table; R5 ("Binding Value Set"
# for Extension.url) happens to duplicate the StructureDefinition URL, and Z6
# holds the ValueSet's Binding Value Set URL.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")
$ws2 = $wb.Worksheets.Item("Elements")

$newSdUrl = "https://2rdoc.pt/ig/ios-lifestyle-medicine/StructureDefinition/social-context"
$newDate  = "2025-08-20T10:40:04+01:00"
$newVsUrl = "https://2rdoc.pt/ig/ios-lifestyle-medicine/ValueSet/social-context-vs"

$ws1.Range("B2").Value = $newSdUrl
$ws2.Range("R5").Value = $newSdUrl

$ws1.Range("B8").Value = $newDate

$ws2.Range("Z6").Value = $newVsUrl

# Re-measured ("bestFit") column widths on the Elements sheet (in Excel
# "characters" units, as used by the ColumnWidth property).
$ws2.Columns.Item(1).ColumnWidth  = 15.666666666666666
$ws2.Columns.Item(2).ColumnWidth  = 15.666666666666666
$ws2.Columns.Item(3).ColumnWidth  = 9.0
$ws2.Columns.Item(4).ColumnWidth  = 6.166666666666667
$ws2.Columns.Item(5).ColumnWidth  = 4.5
$ws2.Columns.Item(6).ColumnWidth  = 3.1666666666666665
$ws2.Columns.Item(7).ColumnWidth  = 3.4999999999999996
$ws2.Columns.Item(8).ColumnWidth  = 11.833333333333332
$ws2.Columns.Item(9).ColumnWidth  = 9.666666666666666
$ws2.Columns.Item(11).ColumnWidth = 13.5
$ws2.Columns.Item(15).ColumnWidth = 11.5
$ws2.Columns.Item(20).ColumnWidth = 7.0
$ws2.Columns.Item(21).ColumnWidth = 12.833333333333332
$ws2.Columns.Item(22).ColumnWidth = 13.166666666666666
$ws2.Columns.Item(23).ColumnWidth = 14.166666666666666
$ws2.Columns.Item(24).ColumnWidth = 13.833333333333332
$ws2.Columns.Item(25).ColumnWidth = 16.166666666666668
$ws2.Columns.Item(26).ColumnWidth = 54.666666666666664
$ws2.Columns.Item(27).ColumnWidth = 4.166666666666667
$ws2.Columns.Item(28).ColumnWidth = 17.166666666666668
$ws2.Columns.Item(29).ColumnWidth = 33.666666666666664
$ws2.Columns.Item(30).ColumnWidth = 12.666666666666666
$ws2.Columns.Item(31).ColumnWidth = 10.5
$ws2.Columns.Item(32).ColumnWidth = 14.166666666666666
$ws2.Columns.Item(33).ColumnWidth = 7.333333333333333
$ws2.Columns.Item(34).ColumnWidth = 7.666666666666667
$ws2.Columns.Item(37).ColumnWidth = 18.666666666666668

# Setting ColumnWidth drops the "hidden" flag on these columns, so restore it.
$ws2.Columns.Item(3).Hidden  = $true
$ws2.Columns.Item(4).Hidden  = $true
$ws2.Columns.Item(31).Hidden = $true
$ws2.Columns.Item(32).Hidden = $true
$ws2.Columns.Item(33).Hidden = $true
